$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.081.90'
$ws.Range('E2').Value = '  -0.06%  '
$ws.Range('D3').Value = '1.621.84'
$ws.Range('E3').Value = '  -0.92%  '
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').Value = '''214.03'
$ws.Range('E5').Value = '  -1.26%  '
$ws.Range('D6').Value = '''0.514'
$ws.Range('E6').Value = '  -0.44%  '
$ws.Range('E7').Value = '  -0.15%  '
$ws.Range('E8').Value = '  +0.52%  '
$ws.Range('D10').Value = '''19.94'
$ws.Range('E10').Value = '  +0.10%  '
$ws.Range('D11').Value = '''0.0841'
$ws.Range('E11').Value = '  -0.73%  '
$ws.Range('D12').Value = '1.848.12'
$ws.Range('E12').Value = '  -0.97%  '
$ws.Range('D13').Value = '1.623.59'
$ws.Range('E13').Value = '  -0.89%  '
$ws.Range('E14').Value = '  +0.00%  '
$ws.Range('D15').Value = '''0.539'
$ws.Range('E15').Value = '  -0.41%  '
$ws.Range('D16').Value = '27.043.01'
$ws.Range('E16').Value = '  -0.18%  '
$ws.Range('D17').Value = '''64.49'
$ws.Range('E17').Value = '  -3.23%  '
$ws.Range('D18').Value = '0.0₃0738'
$ws.Range('E18').Value = '  -0.06%  '
$ws.Range('D19').Value = '''214.16'
$ws.Range('E19').Value = '  -1.22%  '
$ws.Range('E20').Value = '  -0.13%  '
$ws.Range('E21').Value = '  -0.66%  '
$ws.Range('E23').Value = '  -7.59%  '
$ws.Range('E24').Value = '  -0.91%  '
$ws.Range('D25').Value = '''147.74'
$ws.Range('E25').Value = '  +0.71%  '
$ws.Range('B26').Value = 'Cosmos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D26').Value = '''7.42'
$ws.Range('E26').Value = '  +0.29%  '
$ws.Range('B27').Value = 'BinanceUSD'
$ws.Range('C27').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D27').Value = '''1.00'
$ws.Range('E27').Value = '  -0.23%  '
$ws.Range('E28').Value = '  -3.31%  '
$ws.Range('E29').Value = '  -1.07%  '
$ws.Range('E30').Value = '  +0.57%  '
$ws.Range('E31').Value = '  -1.01%  '
$ws.Range('E32').Value = '  -1.72%  '
$ws.Range('D33').Value = '''0.733'
$ws.Range('E33').Value = '  +35.33%  '
$ws.Range('E34').Value = '  -0.42%  '
$ws.Range('D35').Value = '1.338.20'
$ws.Range('E35').Value = '  +2.85%  '
$ws.Range('D36').Value = '''1.56'
$ws.Range('E36').Value = '  -0.77%  '
$ws.Range('E37').Value = '  -0.74%  '
$ws.Range('E38').Value = '  -0.14%  '
$ws.Range('E39').Value = '  -1.59%  '
$ws.Range('E40').Value = '  -0.19%  '
$ws.Range('E41').Value = '  +0.38%  '
$ws.Range('E42').Value = '  -1.63%  '
$ws.Range('D43').Value = '''5.33'
$ws.Range('E43').Value = '  +0.44%  '
$ws.Range('D44').Value = '''63.80'
$ws.Range('E44').Value = '  +3.47%  '
$ws.Range('D45').Value = '1.759.78'
$ws.Range('E45').Value = '  -0.96%  '
$ws.Range('D46').Value = '''89.85'
$ws.Range('E46').Value = '  -1.44%  '
$ws.Range('D47').Value = '''1.64'
$ws.Range('E47').Value = '  +2.60%  '
$ws.Range('D48').Value = '''0.851'
$ws.Range('E48').Value = '  +27.38%  '
$ws.Range('E49').Value = '  +4.71%  '
$ws.Range('E50').Value = '  +0.04%  '
$ws.Range('D51').Value = '''7.58'
$ws.Range('E51').Value = '  -0.78%  '
